# Apply the "add the comparison of CitySim. minor fix and rename column archive"
# edit to the workbook.
#
# 1) Rename the model headers in row 1 (B1:I1) by appending "(APC)".
# 2) Correct a handful of distance_% values in column M with more precise figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename headers (append "(APC)") -----------------------------------
$headerCells = @("B1", "C1", "D1", "E1", "F1", "G1", "H1", "I1")
foreach ($addr in $headerCells) {
    $cell = $ws.Range($addr)
    $oldValue = $cell.Value()
    $cell.Value = $oldValue + "(APC)"
}

# --- 2) Fix distance_% values in column M ----------------------------------
$mFixes = @{
    "M3"  = -1.2
    "M18" = -16.3
    "M19" = -0.4
    "M20" = -2.8
    "M21" = -1.3
    "M24" = -3.2
    "M30" = -2.5
}
foreach ($addr in $mFixes.Keys) {
    $ws.Range($addr).Value2 = $mFixes[$addr]
}

$wb.Save()
